# Updates cryptos list cell values (price + volume%) and restores the
# original row order for the six coins that were re-ranked.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.102.65"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "3.513.92"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'602.43"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").Value = "'183.89"
$ws.Range("E6").Value = "  +5.94%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").Value = "'0.140"
$ws.Range("E9").Value = "  +4.20%  "
$ws.Range("D10").Value = "'7.14"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").Value = "'0.437"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "4.114.02"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "'32.44"
$ws.Range("E13").Value = "  +12.83%  "
$ws.Range("D14").Value = "'0.135"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "68.034.68"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "'0.0000183"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").Value = "3.511.22"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "'6.41"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "'14.80"
$ws.Range("E19").Value = "  +4.01%  "
$ws.Range("D20").Value = "'398.22"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "'8.11"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "'0.547"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'73.38"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000126"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'5.70"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'10.74"
$ws.Range("E27").Value = "  +5.60%  "
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'6.32"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "'1.48"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").Value = "'2.08"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "'24.20"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").Value = "'7.48"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "'1.69"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("D37").Value = "'164.40"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "'1.98"
$ws.Range("E38").Value = "  +3.55%  "
$ws.Range("D39").Value = "'0.877"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").Value = "'7.18"
$ws.Range("E40").Value = "  +4.10%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.78"
$ws.Range("E41").Value = "  +6.64%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.78"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("D43").Value = "'27.99"
$ws.Range("E43").Value = "  +3.48%  "
$ws.Range("D44").Value = "'26.80"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").Value = "'0.0742"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "2.847.87"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").Value = "'42.38"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "'0.0307"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'347.11"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "'34.10"
$ws.Range("E51").Value = "  +2.16%  "
